$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)

$shp.TextFrame.TextRange.Text = "https://github.com/jyothika19/steganography_project.git"
$shp.TextFrame.TextRange.Font.Size = 32
$shp.TextFrame.WordWrap = -1

$shp.Left = 56.11764144897461
$shp.Top = 180
$shp.Width = 800.4706299212598
$shp.Height = 46.04527559055118
